$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original row 13 (A13 empty, B13/C13 = "984972 - Hugo Ricardo Zschommler Sandim")
# is removed entirely; every row below it shifts up by one.
$ws.Rows.Item(13).Delete()

# --- Content replacements (post-shift row numbers) ---

# Row 10 (Objetivos:): body text replaced with the responsible professor string.
$ws.Range("B10").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C10").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 13 (Programa resumido:): body text replaced with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:): body text replaced with the activation date string.
# (leading apostrophe forces text so Excel doesn't coerce it to a date serial)
$ws.Range("B15").Value = "'01/01/2012"
$ws.Range("C15").Value = "'01/01/2012"

# Row 18 (Método:): body text replaced with the responsible professor string again.
$ws.Range("B18").Value = "984972 - Hugo Ricardo Zschommler Sandim"
$ws.Range("C18").Value = "984972 - Hugo Ricardo Zschommler Sandim"

# Row 19 (Critério:): body text becomes the old "Método:" teaching description.
$ws.Range("B19").Value = "Aulas expositivas e aulas práticas de demonstração em oficina. Visita a feiras."
$ws.Range("C19").Value = "Aulas expositivas e aulas práticas de demonstração em oficina. Visita a feiras."

# Row 20 (Norma de recuperação:): body text becomes the old "Critério:" grading formula.
$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

# Row 21 (Bibliografia:): body text becomes the old "Norma de recuperação:" description.
$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
